$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each year block of 4 rows (A/B/C/D sub-periods), the "B" and "C"
# sub-rows swap places (their A:E content swaps), while the "A" and "D"
# sub-rows stay put. Stage the swap through a scratch row via Copy so the
# cell types (numbers / blank cells) round-trip faithfully instead of
# going through .Value assignment. Each destination is explicitly cleared
# first since Copy()/Cut() here won't blank out a destination cell whose
# source is itself blank.
$scratchRow = 1000

for ($start = 2; $start -le 78; $start += 4) {
    $rowB = $start + 1
    $rowC = $start + 2

    $scratchRange = $ws.Range("A" + $scratchRow + ":E" + $scratchRow)
    $rangeB = $ws.Range("A" + $rowB + ":E" + $rowB)
    $rangeC = $ws.Range("A" + $rowC + ":E" + $rowC)

    $rangeB.Copy($scratchRange)
    $rangeB.Clear()
    $rangeC.Copy($rangeB)
    $rangeC.Clear()
    $scratchRange.Copy($rangeC)
    $scratchRange.Clear()
}

# Drop the F/G columns (产销率 / 销售量 non-cumulative figures) entirely.
$ws.Range("F1:G81").EntireColumn.Delete()
